# Dashboard page changes: rename Registration sheet -> TestData, refresh its
# contents with login-page test data (header row + one data row), and tidy
# up the Config sheet's A3 cell formatting. Finish with TestData active.

$wb = $excel.ActiveWorkbook
$wsConfig = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)

# --- Config sheet: A3 no longer needs the quote-prefix style ---
$wsConfig.Cells.Item(3, 1).Style = "Normal"

# --- Rename the Registration sheet ---
$wsData.Name = "TestData"

# --- Reset old content/formatting so new cells start from a clean slate ---
$wsData.Cells.Item(2, 1).Style = "Normal"

# --- New header row ---
$wsData.Cells.Item(1, 1).Value = "TestCaseName"
$wsData.Cells.Item(1, 2).Value = "UserName"
$wsData.Cells.Item(1, 3).Value = "Password"

# --- New data row (content first, formatting applied after) ---
$wsData.Cells.Item(2, 1).Value = "LoginPageTest"
$wsData.Cells.Item(2, 2).Value = "madhur_b"
$wsData.Cells.Item(2, 3).Value = "madhurb"

# --- Header formatting: reuse the bold/orange-fill style already used on Config!A1:C1 ---
$wsConfig.Range("A1:C1").Copy()
$wsData.Range("A1:C1").PasteSpecial(-4122)

# --- B2 (UserName value) keeps the quote-prefix style, like Config!B3 ---
$wsConfig.Cells.Item(3, 2).Copy()
$wsData.Cells.Item(2, 2).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column widths for the new C column / refreshed A & B columns ---
$wsData.Columns.Item(1).ColumnWidth = 149 / 12
$wsData.Columns.Item(2).ColumnWidth = 107 / 12
$wsData.Columns.Item(3).ColumnWidth = 97 / 12

# --- Activate TestData, select C2, zoom 170% ---
$wsData.Activate()
$wsData.Range("C2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 170
